# The deck's slide master, every slide layout, and the notes master each
# carry a "datetimeFigureOut" date placeholder whose cached text had gone
# stale (10/18/2025). Refresh them all to 10/19/2025, matching the cleanup
# described in the commit message.

$p = $ppt.ActivePresentation
$oldDate = "10/18/2025"
$newDate = "10/19/2025"

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $phType = $null
        try { $phType = $shp.PlaceholderFormat.Type } catch { $phType = $null }
        if ($phType -eq 16) {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# 1) Slide master
Update-DatePlaceholders $p.SlideMaster.Shapes

# 2) Every slide layout under the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholders $layouts.Item($L).Shapes
}

# 3) Notes master
Update-DatePlaceholders $p.NotesMaster.Shapes
